$p = $ppt.ActivePresentation

# Slide 1: Title + Subtitle
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "{{TITLE_SLIDE_1}}"
$s1.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "{{SUBTITLE_SLIDE_1}}"

# Slides 2-6: Title + Right Content, each following the same pattern
for ($i = 2; $i -le 6; $i++) {
    $s = $p.Slides.Item($i)
    $s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "{{TITLE_SLIDE_$i}}"
    $s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "{{RIGHT_CONTENT_SLIDE_" + $i + "_1}}"
}
